$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as plain text in this sheet
# (e.g. "26.685.70", "1.000", "  +1.05%  "). Force a text format before
# writing so Excel does not reinterpret numeric-looking strings as numbers,
# then restore the default style so formatting matches the original sheet.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.685.70'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').Value = '1.850.55'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '262.53'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.5377'
$ws.Range('E7').Value = '  +3.40%  '
$ws.Range('D8').Value = '0.3202'
$ws.Range('E8').Value = '  -2.06%  '
$ws.Range('D9').Value = '0.07059'
$ws.Range('E9').Value = '  +3.72%  '
$ws.Range('D10').Value = '19.05'
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('D11').Value = '0.7775'
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = '0.07820'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '1.862.48'
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').Value = '89.54'
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('D15').Value = '5.053'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').Value = '14.15'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '0.000008018'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '26.717.46'
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('D21').Value = '2.088.23'
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('D22').Value = '4.653'
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('D23').Value = '6.057'
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('D24').Value = '9.422'
$ws.Range('E24').Value = '  -1.53%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '142.82'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '2.221'
$ws.Range('E26').Value = '  +1.74%  '
$ws.Range('D27').Value = '1.707'
$ws.Range('E27').Value = '  +2.94%  '
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('D29').Value = '111.90'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').Value = '4.297'
$ws.Range('E30').Value = '  +2.89%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.08752'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.119'
$ws.Range('E32').Value = '  -0.41%  '
$ws.Range('D33').Value = '0.04884'
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('D34').Value = '0.7391'
$ws.Range('E34').Value = '  +2.69%  '
$ws.Range('D35').Value = '1.146'
$ws.Range('E35').Value = '  +1.12%  '
$ws.Range('D36').Value = '2.871'
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('D37').Value = '3.117'
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('D38').Value = '2.360'
$ws.Range('E38').Value = '  +6.84%  '
$ws.Range('D39').Value = '0.01755'
$ws.Range('E39').Value = '  -1.35%  '
$ws.Range('D40').Value = '0.4844'
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('D41').Value = '0.9129'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').Value = '109.49'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').Value = '5.921'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D44').Value = '0.9992'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = '7.729'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').Value = '0.4216'
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('D47').Value = '9.127'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').Value = '0.1258'
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('D49').Value = '35.07'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').Value = '0.05841'
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('D51').Value = '0.8988'
$ws.Range('E51').Value = '  +1.18%  '

$ws.Range("D2:E51").Style = "Normal"
